$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new test-case row with sample valid login credentials.
$ws.Range("B2").Value = "TestPass"
$ws.Range("A2").Value = "JohnTester"

# Update header row: A1 was "email", change to "login". Re-affirm B1 "password".
$ws.Range("A1").Value = "login"
$ws.Range("B1").Value = "password"

# Update the active selection to match the new state.
$ws.Range("B9").Select()
